# Updates the cryptocurrency price (D) and 1h volume-change (E) columns
# per the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.618.33"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.12%  "
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.527.46"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.87%  "
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.10%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.82%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.99%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.562"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.33%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.08%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.519"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.59%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.18"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.11%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0800"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.43%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.110"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.49%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.21"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.02%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.915.26"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.79%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.540.74"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.96%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.20"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.54%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.811"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.40%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.628.42"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.04%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.58"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.86%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0939"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.19%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.86%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.33%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "241.86"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.33%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.94%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.98"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.82%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.02%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.52"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.25%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.25"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.42%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.81%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.58"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.11%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.91"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.22%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "155.64"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.28%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.72"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.43%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.64"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.93%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0782"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.00%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.13"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.91%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.97"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.47%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.52"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.45%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.108"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.67%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.117"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.89%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.24"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.01%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.67"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.79%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.25%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.025.62"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.49%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0295"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.48%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.21"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.13%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.59%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.767.71"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.89%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "79.86"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.18%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "71.78"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.59%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.187"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.58%  "
